# Fruta / hortaliza, semanal
# Insert a new weekly price row for "Vega Modelo de Temuco - Mango" at row 397,
# pushing the existing rows 397-492 down to 398-493.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 397 (shifts 397..492 -> 398..493)
$ws.Rows.Item(397).Insert()

# Populate the newly inserted row with the new week's record
$ws.Range("A397").Value = 10
$ws.Range("B397").Value = "Vega Modelo de Temuco"
$ws.Range("C397").Value = "La Araucanía"
$ws.Range("D397").Value = 44943
$ws.Range("E397").Value = 9
$ws.Range("F397").Value = "Fruta"
$ws.Range("G397").Value = 100108
$ws.Range("H397").Value = "Tropicales y subtropicales"
$ws.Range("I397").Value = 100108002
$ws.Range("J397").Value = "Mango"
$ws.Range("K397").Value = "Sin especificar"
$ws.Range("L397").Value = "Primera"
$ws.Range("M397").Value = 310
$ws.Range("N397").Value = 7500
$ws.Range("O397").Value = 8000
$ws.Range("P397").Value = 7702
$ws.Range("Q397").Value = "$/bandeja 4 kilos"
$ws.Range("R397").Value = "Perú"
$ws.Range("S397").Value = 1926
$ws.Range("T397").Value = 4
